# "Code added for Manage Requests"
# Build out the SearchTestData sheet with the four "accept/withdraw/complete/decline"
# search-skill columns that back the new Manage Requests automation, and move the
# active tab/selection to that sheet the way the author last left the workbook.

$wb = $excel.ActiveWorkbook

$manageRequests = $wb.Worksheets.Item("ManageRequestsTestData")
$search         = $wb.Worksheets.Item("SearchTestData")

# --- SearchTestData: replace the old single "SearchSkill" column with four
#     skill-search columns (Accept / Decline / Withdraw / Complete). Values are
#     entered in the same left-to-right, header-then-data order the workbook
#     author used so new shared-string entries land in a matching order.
$search.Range("A2").Value = "Skill3"
$search.Range("B2").Value = "Skill4"
$search.Range("A1").Value = "SearchSkillToAccept"
$search.Range("C1").Value = "SearchSkillToWithdraw"
$search.Range("D1").Value = "SearchSkillToComplete"
$search.Range("B1").Value = "SearchSkillToDecline"
$search.Range("C2").Value = "Skill5"
$search.Range("D2").Value = "Skill6"

# Match the header styling used elsewhere in the workbook (e.g. ManageRequestsTestData)
$manageRequests.Range("A1:D1").Copy() | Out-Null
$search.Range("A1:D1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Widen the columns to fit the new (longer) header text.
$search.Columns("A").ColumnWidth = 19.998697916666668
$search.Columns("B").ColumnWidth = 20.666666666666668
$search.Columns("C").ColumnWidth = 18.998697916666668
$search.Columns("D").ColumnWidth = 18.666666666666668

$search.PageSetup.PaperSize = [Microsoft.Office.Interop.Excel.XlPaperSize]::xlPaperA4
$search.PageSetup.PaperSize = 9
$search.PageSetup.Orientation = [Microsoft.Office.Interop.Excel.XlPageOrientation]::xlPortrait

# --- ManageRequestsTestData: selection moves to the header row.
$manageRequests.Range("A1:D1").Select() | Out-Null

# --- SearchTestData becomes the active tab (and its selection moves onto the
#     new data), matching where the author left the workbook on save.
$search.Activate()
$search.Range("C13").Select() | Out-Null
